# Data refresh of the "BASE DATA (wajib update)" monitoring sheet.
# - ACTUAL END (col K) moves from 45900 (2025-08-27) to 45879 (2025-08-06) for every task row.
# - % COMPLETE (col L) is updated for the rows whose progress changed since the last pull.
# - The active selection / scroll position on the sheet is refreshed to reflect where the
#   author left off (single cell M93 instead of the old H2:M127 block selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BASE DATA (wajib update)")

# --- ACTUAL END: every data row (2-127) now reports 45879 instead of 45900 ---
$ws.Range("K2:K127").Value = 45879

# --- % COMPLETE: refreshed progress figures for the rows that moved ---
$pctComplete = @{
    2   = 0.7374
    11  = 0.0478
    20  = 0.1507
    45  = 0.3216
    66  = 0.1716
    69  = 0.9006
    70  = 0.9989
    80  = 0.8049
    84  = 0.8598
    88  = 0.703
    91  = 0.002934
    99  = 0.2936
    102 = 0.338
}

foreach ($row in $pctComplete.Keys) {
    $ws.Range("L" + $row).Value = $pctComplete[$row]
}

# --- Refresh the sheet's active selection to where the author left off ---
$ws.Activate()
$ws.Range("M93").Select()
